$d = $word.ActiveDocument

# Q4: replace full text
$d.Content.Find.Execute(
    "Q4: What did you find easy about development in this framework?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Q4: What appear to be benefits to using web technologies for development on this device? Conversely, what advantages might a native framework have over the framework you used?",
    2
)

# Q5: replace full text
$d.Content.Find.Execute(
    "Q5: What did you find challenging about development in this framework?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Q5: What did you find easy and challenging about development in this framework?",
    2
)
